# Daily attendance processing - 2026-01-12 17:38:25
# Rotate the "Recorded By" (column G) comma-separated list values:
# move the first item of each list to the end (left-rotate by one),
# leaving single-valued cells untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }
    if ($val -notlike "*,*") { continue }

    $parts = $val -split ",\s*"
    if ($parts.Count -lt 2) { continue }

    $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
    $cell.Value2 = $rotated
}
